$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 5146.4375
$ws.Range("I62").Value = 3949.3
$ws.Range("K62").Value = 3949.3
$ws.Range("M62").Value = -3325.3
$ws.Range("H65").Value = 5146.4375
$ws.Range("I65").Value = 3949.3
$ws.Range("K65").Value = 19746.5
$ws.Range("M65").Value = -16626.5
$ws.Range("H113").Value = 17252828
$ws.Range("I113").Value = 5008950
$ws.Range("K113").Value = 5008950
$ws.Range("M113").Value = -5005696
$ws.Range("H137").Value = 6737.5864
$ws.Range("I137").Value = 2358.7693
$ws.Range("K137").Value = 7076.3079
$ws.Range("M137").Value = -4526.3079

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7361036
$ws.Range("I32").Value = 7940192.5
$ws.Range("K32").Value = 7940192.5
$ws.Range("M32").Value = -7939905.5
$ws.Range("H37").Value = 34000
$ws.Range("J37").Value = 34000
$ws.Range("L37").Value = 34000
$ws.Range("N37").Value = -34546
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()
$ws.Range("H49").Value = 0
$ws.Range("J49").Value = 0
$ws.Range("L49").Value = 0
$ws.Range("N49").ClearContents()
$ws.Range("H61").Value = 16708080
$ws.Range("I61").Value = 19237612
$ws.Range("K61").Value = 19237612
$ws.Range("M61").Value = -19237400
$ws.Range("H96").Value = 109995
$ws.Range("J96").Value = 109995
$ws.Range("L96").Value = 109995
$ws.Range("N96").Value = -115487
$ws.Range("H97").Value = 1848.7142
$ws.Range("I97").Value = 1848.7142
$ws.Range("K97").Value = 1848.7142
$ws.Range("M97").Value = -1352.7142
$ws.Range("H102").Value = 2383
$ws.Range("I102").Value = 2383
$ws.Range("K102").Value = 2383
$ws.Range("M102").Value = -761
$ws.Range("H131").Value = 146000
$ws.Range("J131").Value = 146000
$ws.Range("L131").Value = 146000
$ws.Range("N131").Value = -156080
$ws.Range("H136").Value = 16708080
$ws.Range("I136").Value = 19237612
$ws.Range("K136").Value = 57712836
$ws.Range("M136").Value = -57710286

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 146.85715
$ws.Range("I22").Value = 121.5
$ws.Range("K22").Value = 121.5
$ws.Range("M22").Value = 51.5
$ws.Range("H70").Value = 253330
$ws.Range("J70").Value = 253330
$ws.Range("L70").Value = 253330
$ws.Range("N70").Value = -253916
$ws.Range("H73").Value = 253330
$ws.Range("J73").Value = 253330
$ws.Range("L73").Value = 253330
$ws.Range("N73").Value = -255358
$ws.Range("H80").Value = 1431.375
$ws.Range("J80").Value = 2037.2
$ws.Range("L80").Value = 2037.2
$ws.Range("N80").Value = -4033.2
$ws.Range("H83").Value = 1431.375
$ws.Range("J83").Value = 2037.2
$ws.Range("L83").Value = 10186
$ws.Range("N83").Value = -20170
$ws.Range("H86").Value = 2401.5715
$ws.Range("I86").Value = 2401.5715
$ws.Range("K86").Value = 2401.5715
$ws.Range("M86").Value = -1278.5715
$ws.Range("H89").Value = 2401.5715
$ws.Range("I89").Value = 2401.5715
$ws.Range("K89").Value = 12007.8575
$ws.Range("M89").Value = -6391.8575
$ws.Range("H99").Value = 7785.3887
$ws.Range("I99").Value = 9410.538
$ws.Range("K99").Value = 9410.538
$ws.Range("M99").Value = -7912.538

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 2460.05
$ws.Range("I105").Value = 2630
$ws.Range("J105").Value = 2144.4285
$ws.Range("K105").Value = 2630
$ws.Range("L105").Value = 2144.4285
$ws.Range("M105").Value = -883
$ws.Range("N105").Value = -5638.4285
$ws.Range("H132").Value = 4967.75
$ws.Range("I132").Value = 4820.2856
$ws.Range("K132").Value = 14460.8568
$ws.Range("M132").Value = -11930.8568

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H15").Value = 152.3077
$ws.Range("I15").Value = 152.3077
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 456.9231
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -316.9231
$ws.Range("N15").ClearContents()
$ws.Range("H55").Value = 7692.364
$ws.Range("I55").Value = 7086.6
$ws.Range("J55").Value = 8197.166999999999
$ws.Range("K55").Value = 21259.8
$ws.Range("L55").Value = 24591.501
$ws.Range("M55").Value = -21082.8
$ws.Range("N55").Value = -24945.501
$ws.Range("H64").Value = 5050
$ws.Range("I64").Value = 5050
$ws.Range("K64").Value = 15150
$ws.Range("M64").Value = -14880
$ws.Range("H67").Value = 5050
$ws.Range("I67").Value = 5050
$ws.Range("K67").Value = 15150
$ws.Range("M67").Value = -14214
$ws.Range("H107").Value = 665.7143
$ws.Range("I107").Value = 487.4
$ws.Range("K107").Value = 1462.2
$ws.Range("M107").Value = 457.8000000000002
$ws.Range("H125").Value = 15005.5
$ws.Range("J125").Value = 15005.5
$ws.Range("L125").Value = 45016.5
$ws.Range("N125").Value = -54856.5

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 158.2
$ws.Range("I2").Value = 197
$ws.Range("K2").Value = 197
$ws.Range("M2").Value = -84
$ws.Range("H35").Value = 92203.75
$ws.Range("I35").Value = 106605
$ws.Range("K35").Value = 106605
$ws.Range("M35").Value = -106307
$ws.Range("H53").Value = 0
$ws.Range("J53").Value = 0
$ws.Range("L53").Value = 0
$ws.Range("N53").ClearContents()
$ws.Range("H70").Value = 4500
$ws.Range("I70").Value = 500
$ws.Range("J70").Value = 6500
$ws.Range("K70").Value = 500
$ws.Range("L70").Value = 6500
$ws.Range("M70").Value = -230
$ws.Range("N70").Value = -7040
$ws.Range("H73").Value = 4500
$ws.Range("I73").Value = 500
$ws.Range("J73").Value = 6500
$ws.Range("K73").Value = 500
$ws.Range("L73").Value = 6500
$ws.Range("M73").Value = 436
$ws.Range("N73").Value = -8372
$ws.Range("H80").Value = 1282.5
$ws.Range("I80").Value = 1465
$ws.Range("J80").Value = 1100
$ws.Range("K80").Value = 1465
$ws.Range("L80").Value = 1100
$ws.Range("M80").Value = -467
$ws.Range("N80").Value = -3096
$ws.Range("H83").Value = 1282.5
$ws.Range("I83").Value = 1465
$ws.Range("J83").Value = 1100
$ws.Range("K83").Value = 7325
$ws.Range("L83").Value = 5500
$ws.Range("M83").Value = -2333
$ws.Range("N83").Value = -15484
$ws.Range("H102").Value = 2488.75
$ws.Range("I102").Value = 1833.3846
$ws.Range("K102").Value = 1833.3846
$ws.Range("M102").Value = -211.3846000000001

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H29").Value = 16000
$ws.Range("I29").Value = 0
$ws.Range("J29").Value = 16000
$ws.Range("K29").Value = 0
$ws.Range("L29").Value = 16000
$ws.Range("N29").Value = -16590
$ws.Range("M29").ClearContents()
$ws.Range("H33").Value = 0
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()
$ws.Range("H34").Value = 36500
$ws.Range("I34").Value = 36500
$ws.Range("K34").Value = 36500
$ws.Range("M34").Value = -36328
$ws.Range("H57").Value = 19888
$ws.Range("J57").Value = 19888
$ws.Range("L57").Value = 19888
$ws.Range("N57").Value = -21020
$ws.Range("H68").Value = 3700
$ws.Range("I68").Value = 3500
$ws.Range("K68").Value = 3500
$ws.Range("M68").Value = -2751
$ws.Range("H71").Value = 3700
$ws.Range("I71").Value = 3500
$ws.Range("K71").Value = 17500
$ws.Range("M71").Value = -13756
$ws.Range("H132").Value = 7340716.5
$ws.Range("I132").Value = 841062.3
$ws.Range("K132").Value = 2523186.9
$ws.Range("M132").Value = -2520656.9

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H120").Value = 105000
$ws.Range("J120").Value = 105000
$ws.Range("L120").Value = 105000
$ws.Range("N120").Value = -114676
$ws.Range("H132").Value = 281322.84
$ws.Range("J132").Value = 1258623.4
$ws.Range("L132").Value = 3775870.2
$ws.Range("N132").Value = -3780930.2
$ws.Range("H136").Value = 2801.85
$ws.Range("I136").Value = 1002.3125
$ws.Range("K136").Value = 3006.9375
$ws.Range("M136").Value = -456.9375
